$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite Sheet1 with the ShiftBookExcel2 header row.
$ws.Range("A1").Value = "timeStamp"
$ws.Range("B1").Value = "studentID"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "date"
$ws.Range("E1").Value = "shift"
